# 🔄 Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resultados (columna G) y profit (columna H) para filas recién resueltas
$ws.Range("G98").Value = "Fallo"
$ws.Range("H98").Value = -1

$ws.Range("G100").Value = "Acierto"
$ws.Range("H100").Value = 1.75

$ws.Range("G131").Value = "Fallo"
$ws.Range("H131").Value = -1

$ws.Range("G133").Value = "Acierto"
$ws.Range("H133").Value = 4

$ws.Range("G134").Value = "Fallo"
$ws.Range("H134").Value = -1

$ws.Range("G139").Value = "Fallo"
$ws.Range("H139").Value = -1

$ws.Range("G140").Value = "Fallo"
$ws.Range("H140").Value = -1

$ws.Range("G141").Value = "Acierto"
$ws.Range("H141").Value = 1.5

$ws.Range("G144").Value = "Acierto"
$ws.Range("H144").Value = 2

$ws.Range("G145").Value = "Fallo"
$ws.Range("H145").Value = -1

$ws.Range("G149").Value = "Fallo"
$ws.Range("H149").Value = -1

# event_id de las filas 152 y 153 pasan de texto a numerico
$ws.Range("A152").Value = 14552522
$ws.Range("A153").Value = 14552653

Write-Host "Tracker actualizado"
